$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness) holds a constant 7573 for every logged generation
# (rows 2-252). This run's final fitness value has been corrected to 7310,
# so overwrite the whole column with the new value.
$ws.Range("C2:C252").Value = 7310
